# Lesson Plan for Venv for Robotics Club - 311_dict.xlsx edits
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. D9 "Dept" comment - reflow the line breaks and drop the URL
#    (two runs: bold header run + normal body run)
# ---------------------------------------------------------------
$d9Bold = "There are Nine Departments:`n "
$d9Body = "'Animal Care Services', 'Office of Historic Preservation',  'Public Works', '311', 'Development Services',  `n'Parks and Recreation', 'Human Services', 'Solid Waste Management', 'Metro Health' "
$ws.Range("D9").Value2 = $d9Bold + $d9Body
$ws.Range("D9").Characters(1, $d9Bold.Length).Font.Bold = $true
$ws.Range("D9").Characters($d9Bold.Length + 1, $d9Body.Length).Font.Bold = $false

# ---------------------------------------------------------------
# 2. XCOORD / YCOORD comment - append the NOAA converter link text
#    (D16 and D17 share the same comment text)
# ---------------------------------------------------------------
$nadText = "NAD(1983)State Plane Texas South Central FIPS 4204 Feet `nLink to NOAA Coordinate Converter`nhttps://www.ngs.noaa.gov/NCAT/"
$ws.Range("D16").Value2 = $nadText
$ws.Range("D17").Value2 = $nadText

# ---------------------------------------------------------------
# 3. Wrap text on D16 / D17 so the new multi-line comment is visible
# ---------------------------------------------------------------
$ws.Range("D16").WrapText = $true
$ws.Range("D17").WrapText = $true

# ---------------------------------------------------------------
# 4. Row heights / column width adjustments
# ---------------------------------------------------------------
$ws.Rows.Item(9).RowHeight = 35.05
$ws.Rows.Item(16).RowHeight = 46.25
$ws.Rows.Item(17).RowHeight = 46.25
$ws.Columns.Item(4).ColumnWidth = 95.57

# ---------------------------------------------------------------
# 5. Hyperlinks - drop the old D9 department-directory link, keep the
#    council district link (now the first hyperlink) and add a new
#    hyperlink on D17 pointing at the NOAA coordinate converter.
# ---------------------------------------------------------------
$toRemove = @()
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$9') {
        $toRemove += $h
    }
}
foreach ($h in $toRemove) {
    $h.Delete()
}

$ws.Hyperlinks.Add($ws.Range("D17"), "https://www.ngs.noaa.gov/NCAT/", [Type]::Missing, [Type]::Missing, "https://www.ngs.noaa.gov/NCAT/") | Out-Null

# ---------------------------------------------------------------
# 6. Update the active selection to C18
# ---------------------------------------------------------------
$ws.Range("C18").Select() | Out-Null
